# Apply numeric updates to the Excalibur_Profits sheets (scheduled runner refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5296
$ws.Range("I43").Value = 5694.5
$ws.Range("K43").Value = 5694.5
$ws.Range("M43").Value = -5625.5
$ws.Range("H100").Value = 3744.1333
$ws.Range("I100").Value = 1467.6364
$ws.Range("K100").Value = 1467.6364
$ws.Range("M100").Value = -926.6364000000001
$ws.Range("H111").Value = 4101.6
$ws.Range("J111").Value = 3999.5
$ws.Range("L111").Value = 11998.5
$ws.Range("N111").Value = -18132.5
$ws.Range("H129").Value = 1848.2
$ws.Range("I129").Value = 1424.8462
$ws.Range("J129").Value = 4600
$ws.Range("K129").Value = 4274.5386
$ws.Range("L129").Value = 13800
$ws.Range("M129").Value = 725.4614000000001
$ws.Range("N129").Value = -23800
$ws.Range("H132").Value = 40468.297
$ws.Range("I132").Value = 45624.62
$ws.Range("K132").Value = 136873.86
$ws.Range("M132").Value = -134343.86
$ws.Range("H141").Value = 2109.7778
$ws.Range("I141").Value = 2418.4
$ws.Range("J141").Value = 1724
$ws.Range("K141").Value = 7255.200000000001
$ws.Range("L141").Value = 5172
$ws.Range("M141").Value = -2075.200000000001
$ws.Range("N141").Value = -15532

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1520.7441
$ws.Range("I2").Value = 1449.8918
$ws.Range("J2").Value = 1957.6666
$ws.Range("K2").Value = 1449.8918
$ws.Range("L2").Value = 1957.6666
$ws.Range("M2").Value = -1336.8918
$ws.Range("N2").Value = -2183.6666
$ws.Range("H45").Value = 11842.429
$ws.Range("I45").Value = 9379.6
$ws.Range("J45").Value = 17999.5
$ws.Range("K45").Value = 9379.6
$ws.Range("L45").Value = 17999.5
$ws.Range("M45").Value = -9002.6
$ws.Range("N45").Value = -18753.5
$ws.Range("H74").Value = 2623.7678
$ws.Range("I74").Value = 2067.913
$ws.Range("K74").Value = 2067.913
$ws.Range("M74").Value = -1193.913
$ws.Range("H77").Value = 2623.7678
$ws.Range("I77").Value = 2067.913
$ws.Range("K77").Value = 10339.565
$ws.Range("M77").Value = -5971.565000000001
$ws.Range("H110").Value = 917.94446
$ws.Range("J110").Value = 420
$ws.Range("L110").Value = 420
$ws.Range("N110").Value = -4510
$ws.Range("H116").Value = 1520.7441
$ws.Range("I116").Value = 1449.8918
$ws.Range("J116").Value = 1957.6666
$ws.Range("K116").Value = 1449.8918
$ws.Range("L116").Value = 1957.6666
$ws.Range("M116").Value = 844.1081999999999
$ws.Range("N116").Value = -6545.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1520.7441
$ws.Range("I3").Value = 1449.8918
$ws.Range("J3").Value = 1957.6666
$ws.Range("K3").Value = 1449.8918
$ws.Range("L3").Value = 1957.6666
$ws.Range("M3").Value = -1335.8918
$ws.Range("N3").Value = -2185.6666
$ws.Range("H20").Value = 4412.3
$ws.Range("I20").Value = 3874.8572
$ws.Range("K20").Value = 3874.8572
$ws.Range("M20").Value = -3627.8572
$ws.Range("H86").Value = 2410.389
$ws.Range("I86").Value = 1420.4445
$ws.Range("J86").Value = 3400.3333
$ws.Range("K86").Value = 1420.4445
$ws.Range("L86").Value = 3400.3333
$ws.Range("M86").Value = -297.4445000000001
$ws.Range("N86").Value = -5646.3333
$ws.Range("H89").Value = 2410.389
$ws.Range("I89").Value = 1420.4445
$ws.Range("J89").Value = 3400.3333
$ws.Range("K89").Value = 7102.2225
$ws.Range("L89").Value = 17001.6665
$ws.Range("M89").Value = -1486.2225
$ws.Range("N89").Value = -28233.6665
$ws.Range("H106").Value = 24633
$ws.Range("J106").Value = 24633
$ws.Range("L106").Value = 24633
$ws.Range("N106").Value = -27157
$ws.Range("H107").Value = 1894.3489
$ws.Range("I107").Value = 1989.3158
$ws.Range("J107").Value = 1172.6
$ws.Range("K107").Value = 1989.3158
$ws.Range("L107").Value = 1172.6
$ws.Range("M107").Value = -69.31580000000008
$ws.Range("N107").Value = -5012.6
$ws.Range("H134").Value = 34253.56
$ws.Range("I134").Value = 39622.75
$ws.Range("J134").Value = 9197.333000000001
$ws.Range("K134").Value = 118868.25
$ws.Range("L134").Value = 27591.999
$ws.Range("M134").Value = -116333.25
$ws.Range("N134").Value = -32661.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 67598.25
$ws.Range("J70").Value = 67598.25
$ws.Range("L70").Value = 67598.25
$ws.Range("N70").Value = -68228.25
$ws.Range("H73").Value = 67598.25
$ws.Range("J73").Value = 67598.25
$ws.Range("L73").Value = 67598.25
$ws.Range("N73").Value = -69782.25
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("H105").Value = 1910.2
$ws.Range("I105").Value = 1910.2
$ws.Range("K105").Value = 1910.2
$ws.Range("M105").Value = -163.2
$ws.Range("H134").Value = 2120.8262
$ws.Range("I134").Value = 1376.3182
$ws.Range("K134").Value = 4128.9546
$ws.Range("M134").Value = -1593.9546
$ws.Range("N80").ClearContents()
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 6624.375
$ws.Range("J39").Value = 6624.375
$ws.Range("L39").Value = 19873.125
$ws.Range("N39").Value = -20461.125
$ws.Range("H70").Value = 4664.5
$ws.Range("I70").Value = 2999.25
$ws.Range("K70").Value = 8997.75
$ws.Range("M70").Value = -8682.75
$ws.Range("H73").Value = 4664.5
$ws.Range("I73").Value = 2999.25
$ws.Range("K73").Value = 8997.75
$ws.Range("M73").Value = -7905.75
$ws.Range("H82").Value = 5583.25
$ws.Range("I82").Value = 4121.3335
$ws.Range("J82").Value = 9969
$ws.Range("K82").Value = 12364.0005
$ws.Range("L82").Value = 29907
$ws.Range("M82").Value = -11958.0005
$ws.Range("N82").Value = -30719
$ws.Range("H85").Value = 5583.25
$ws.Range("I85").Value = 4121.3335
$ws.Range("J85").Value = 9969
$ws.Range("K85").Value = 12364.0005
$ws.Range("L85").Value = 29907
$ws.Range("M85").Value = -10960.0005
$ws.Range("N85").Value = -32715
$ws.Range("H102").Value = 7139.25
$ws.Range("I102").Value = 3737.4
$ws.Range("K102").Value = 11212.2
$ws.Range("M102").Value = -8778.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 45250
$ws.Range("J15").Value = 45250
$ws.Range("L15").Value = 45250
$ws.Range("N15").Value = -45826
$ws.Range("H81").Value = 45250
$ws.Range("J81").Value = 45250
$ws.Range("L81").Value = 45250
$ws.Range("N81").Value = -47246
$ws.Range("H84").Value = 45250
$ws.Range("J84").Value = 45250
$ws.Range("L84").Value = 135750
$ws.Range("N84").Value = -145734
$ws.Range("H93").Value = 48000
$ws.Range("J93").Value = 48000
$ws.Range("L93").Value = 48000
$ws.Range("N93").Value = -51744
$ws.Range("H102").Value = 1963.0197
$ws.Range("I102").Value = 1208.6923
$ws.Range("J102").Value = 4414.5835
$ws.Range("K102").Value = 1208.6923
$ws.Range("L102").Value = 4414.5835
$ws.Range("M102").Value = 413.3077000000001
$ws.Range("N102").Value = -7658.5835

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1303.8
$ws.Range("I16").Value = 1142.5
$ws.Range("J16").Value = 1949
$ws.Range("K16").Value = 1142.5
$ws.Range("L16").Value = 1949
$ws.Range("M16").Value = -972.5
$ws.Range("H40").Value = 2161.875
$ws.Range("I40").Value = 2161.875
$ws.Range("K40").Value = 2161.875
$ws.Range("M40").Value = -2025.875
$ws.Range("H46").Value = 3664.5833
$ws.Range("I46").Value = 625
$ws.Range("J46").Value = 3940.9092
$ws.Range("K46").Value = 625
$ws.Range("L46").Value = 3940.9092
$ws.Range("M46").Value = -437
$ws.Range("N46").Value = -4316.9092
$ws.Range("H80").Value = 60051.2
$ws.Range("J80").Value = 60051.2
$ws.Range("L80").Value = 60051.2
$ws.Range("N80").Value = -62297.2
$ws.Range("H82").Value = 2870.8635
$ws.Range("J82").Value = 3680.4614
$ws.Range("L82").Value = 3680.4614
$ws.Range("N82").Value = -4402.4614
$ws.Range("H83").Value = 60051.2
$ws.Range("J83").Value = 60051.2
$ws.Range("L83").Value = 180153.6
$ws.Range("N83").Value = -191385.6
$ws.Range("H85").Value = 2870.8635
$ws.Range("J85").Value = 3680.4614
$ws.Range("L85").Value = 3680.4614
$ws.Range("N85").Value = -6176.4614
$ws.Range("H100").Value = 11477.846
$ws.Range("I100").Value = 2900
$ws.Range("K100").Value = 2900
$ws.Range("M100").Value = -2359
$ws.Range("N16").Value = -2289

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 142868960
$ws.Range("I2").Value = 15723.5
$ws.Range("J2").Value = 333339940
$ws.Range("K2").Value = 15723.5
$ws.Range("L2").Value = 333339940
$ws.Range("M2").Value = -15611.5
$ws.Range("N2").Value = -333340164
$ws.Range("H37").Value = 26499.5
$ws.Range("I37").Value = 20000
$ws.Range("K37").Value = 20000
$ws.Range("H107").Value = 1047.8975
$ws.Range("J107").Value = 1507.6428
$ws.Range("L107").Value = 4522.928400000001
$ws.Range("N107").Value = -8362.928400000001
$ws.Range("M37").Value = -19797
